$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 2P": update Blancos/Reprobados/Aprobados/Por_Apro and add Promedio ---
$ws2p = $wb.Worksheets.Item("Estadisticos 2P")

$ws2p.Range("D2").Value = 0
$ws2p.Range("E2").Value = 0
$ws2p.Range("F2").Value = 41
$ws2p.Range("G2").Value = 100
$ws2p.Range("H2").Value = 8.9

$ws2p.Range("D3").Value = 0
$ws2p.Range("E3").Value = 0
$ws2p.Range("F3").Value = 35
$ws2p.Range("G3").Value = 100
$ws2p.Range("H3").Value = 9

$ws2p.Range("D4").Value = 0
$ws2p.Range("E4").Value = 0
$ws2p.Range("F4").Value = 39
$ws2p.Range("G4").Value = 100
$ws2p.Range("H4").Value = 8.2

$ws2p.Range("D5").Value = 0
$ws2p.Range("E5").Value = 0
$ws2p.Range("F5").Value = 36
$ws2p.Range("G5").Value = 100
$ws2p.Range("H5").Value = 8.5

# --- Sheet "Estadisticos Final": recompute Reprobados/Aprobados/Por_Apro/Promedio ---
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")

$wsFinal.Range("E2").Value = 0
$wsFinal.Range("F2").Value = 41
$wsFinal.Range("G2").Value = 100
$wsFinal.Range("H2").Value = 9.5

$wsFinal.Range("E3").Value = 0
$wsFinal.Range("F3").Value = 35
$wsFinal.Range("G3").Value = 100
$wsFinal.Range("H3").Value = 9.3

$wsFinal.Range("E4").Value = 0
$wsFinal.Range("F4").Value = 39
$wsFinal.Range("G4").Value = 100
$wsFinal.Range("H4").Value = 8.8

$wsFinal.Range("E5").Value = 0
$wsFinal.Range("F5").Value = 36
$wsFinal.Range("G5").Value = 100
$wsFinal.Range("H5").Value = 9.4

# --- Sheet "Rescatables": the rescue list is now empty, remove the student rows ---
$wsResc = $wb.Worksheets.Item("Rescatables")
$wsResc.Range("A2:G5").EntireRow.Delete()
